$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 300.75
$ws.Range("I2").Value = 266.9
$ws.Range("J2").Value = 470
$ws.Range("K2").Value = 266.9
$ws.Range("L2").Value = 470
$ws.Range("M2").Value = -153.9
$ws.Range("N2").Value = -696

$ws.Range("H29").Value = 575
$ws.Range("J29").Value = 975
$ws.Range("L29").Value = 2925
$ws.Range("N29").Value = -3487

$ws.Range("H101").Value = 2084237.5
$ws.Range("I101").Value = 392.66666
$ws.Range("J101").Value = 4168082.2
$ws.Range("K101").Value = 1177.99998
$ws.Range("L101").Value = 12504246.6
$ws.Range("M101").Value = 444.0000199999999
$ws.Range("N101").Value = -12507490.6

$ws.Range("H112").Value = 2931.9473
$ws.Range("J112").Value = 3475.1333
$ws.Range("L112").Value = 10425.3999
$ws.Range("N112").Value = -12641.3999

$ws.Range("H127").Value = 5282
$ws.Range("I127").Value = 6376
$ws.Range("K127").Value = 19128
$ws.Range("M127").Value = -14168

$ws.Range("H129").Value = 1804.9166
$ws.Range("I129").Value = 398.33334
$ws.Range("J129").Value = 2273.7778
$ws.Range("K129").Value = 1195.00002
$ws.Range("L129").Value = 6821.3334
$ws.Range("M129").Value = 3804.99998
$ws.Range("N129").Value = -16821.3334

$ws.Range("H137").Value = 54457.633
$ws.Range("I137").Value = 1595.75
$ws.Range("K137").Value = 4787.25
$ws.Range("M137").Value = -2237.25

$ws.Range("H138").Value = 4013.3447
$ws.Range("J138").Value = 3502.4878
$ws.Range("L138").Value = 10507.4634
$ws.Range("N138").Value = -20787.4634

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 26788.1
$ws.Range("I32").Value = 18012
$ws.Range("J32").Value = 27763.223
$ws.Range("K32").Value = 18012
$ws.Range("L32").Value = 27763.223
$ws.Range("M32").Value = -17725
$ws.Range("N32").Value = -28337.223

$ws.Range("H45").Value = 1580.84
$ws.Range("I45").Value = 1524.0769
$ws.Range("K45").Value = 1524.0769
$ws.Range("M45").Value = -1147.0769

$ws.Range("H61").Value = 22245.592
$ws.Range("I61").Value = 28753.2
$ws.Range("K61").Value = 28753.2
$ws.Range("M61").Value = -28541.2

$ws.Range("H132").Value = 1539.1918
$ws.Range("I132").Value = 1099.7843
$ws.Range("K132").Value = 3299.3529
$ws.Range("M132").Value = -769.3528999999999

$ws.Range("H136").Value = 22245.592
$ws.Range("I136").Value = 28753.2
$ws.Range("K136").Value = 86259.60000000001
$ws.Range("M136").Value = -83709.60000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1312.5714
$ws.Range("I99").Value = 1186.8235
$ws.Range("J99").Value = 1847
$ws.Range("K99").Value = 1186.8235
$ws.Range("L99").Value = 1847
$ws.Range("M99").Value = 311.1765
$ws.Range("N99").Value = -4843

$ws.Range("H125").Value = 29000
$ws.Range("J125").Value = 29000
$ws.Range("L125").Value = 29000
$ws.Range("N125").Value = -38840

$ws.Range("H134").Value = 4465.8696
$ws.Range("I134").Value = 4286.1904
$ws.Range("K134").Value = 12858.5712
$ws.Range("M134").Value = -10323.5712

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2114.697
$ws.Range("I31").Value = 1478.5
$ws.Range("K31").Value = 1478.5
$ws.Range("M31").Value = -1183.5

$ws.Range("H34").Value = 2114.697
$ws.Range("I34").Value = 1478.5
$ws.Range("K34").Value = 1478.5
$ws.Range("M34").Value = -1276.5

$ws.Range("H99").Value = 4728.125
$ws.Range("J99").Value = 6000
$ws.Range("L99").Value = 6000
$ws.Range("N99").Value = -8996

$ws.Range("H126").Value = 4728.125
$ws.Range("J126").Value = 6000
$ws.Range("L126").Value = 18000
$ws.Range("N126").Value = -22940

$ws.Range("H132").Value = 1978.862
$ws.Range("I132").Value = 1338.75
$ws.Range("K132").Value = 4016.25
$ws.Range("M132").Value = -1486.25

$ws.Range("H141").Value = 88734.55499999999
$ws.Range("J141").Value = 90826.375
$ws.Range("L141").Value = 90826.375
$ws.Range("N141").Value = -101186.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 26422.637
$ws.Range("J113").Value = 1503.4722
$ws.Range("L113").Value = 4510.4166
$ws.Range("N113").Value = -8850.4166

$ws.Range("H122").Value = 1420.0588
$ws.Range("J122").Value = 1610.8462
$ws.Range("L122").Value = 14497.6158
$ws.Range("N122").Value = -19397.6158

$ws.Range("H131").Value = 17044.715
$ws.Range("J131").Value = 17389.396
$ws.Range("L131").Value = 52168.188
$ws.Range("N131").Value = -62248.188

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H29").Value = 80008
$ws.Range("J29").Value = 80008
$ws.Range("L29").Value = 80008
$ws.Range("N29").Value = -80588

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1298.3334
$ws.Range("I22").Value = 956.7143
$ws.Range("J22").Value = 1417.9
$ws.Range("K22").Value = 956.7143
$ws.Range("L22").Value = 1417.9
$ws.Range("M22").Value = -661.7143
$ws.Range("N22").Value = -2007.9

$ws.Range("H27").Value = 1298.3334
$ws.Range("I27").Value = 956.7143
$ws.Range("J27").Value = 1417.9
$ws.Range("K27").Value = 956.7143
$ws.Range("L27").Value = 1417.9
$ws.Range("M27").Value = -849.7143

$ws.Range("H46").Value = 2370.6875
$ws.Range("I46").Value = 1985
$ws.Range("J46").Value = 3013.5
$ws.Range("K46").Value = 1985
$ws.Range("L46").Value = 3013.5
$ws.Range("M46").Value = -1797
$ws.Range("N46").Value = -3389.5

$ws.Range("H68").Value = 2623.75
$ws.Range("I68").Value = 1998.2
$ws.Range("K68").Value = 1998.2
$ws.Range("M68").Value = -1249.2

$ws.Range("H71").Value = 2623.75
$ws.Range("I71").Value = 1998.2
$ws.Range("K71").Value = 9991
$ws.Range("M71").Value = -6247

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 13962.2
$ws.Range("J29").Value = 13962.2
$ws.Range("L29").Value = 13962.2
$ws.Range("N29").Value = -14542.2

$ws.Range("H107").Value = 886.3
$ws.Range("I107").Value = 762.55554
$ws.Range("K107").Value = 2287.66662
$ws.Range("M107").Value = -367.66662

$ws.Range("H132").Value = 1679.6444
$ws.Range("I132").Value = 1235.1
$ws.Range("J132").Value = 2035.28
$ws.Range("K132").Value = 3705.3
$ws.Range("L132").Value = 6105.84
$ws.Range("M132").Value = -1175.3
$ws.Range("N132").Value = -11165.84
